$d = $word.ActiveDocument

# Update version string in the Date paragraph: v2.0.beta8 -> v2.0.beta9
$d.Content.Find.Execute("v2.0.beta8", $true, $false, $false, $false, $false,
                         $true, 1, $false, "v2.0.beta9", 2)

# Update release date: (2015-12-04) -> (2015-12-16)
$d.Content.Find.Execute("(2015-12-04)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(2015-12-16)", 2)
